# Generate Report for Handback
# Updates the localization-status report after a successful handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Handback DateTime is refreshed for zh-cn / de-de
#  - The stale "handback file is not the latest" error is cleared
#  - A couple of columns are re-sized to fit the new content

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns for zh-cn / de-de ---
$ws_overview.Range("E2").Value = $newStatus
$ws_overview.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$ws_zhcn.Range("C2").Value = $newStatus
$ws_zhcn.Range("K2").Value = "2016-09-01 00:53:40"
$ws_zhcn.Range("P2").Value = ""

# --- de-de sheet ---
$ws_dede.Range("C2").Value = $newStatus
$ws_dede.Range("K2").Value = "2016-09-01 00:53:47"
$ws_dede.Range("P2").Value = ""

# --- Column width adjustments (content-driven resize) ---
# Overview: zh-cn / de-de status columns widened
$ws_overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$ws_overview.Columns.Item(6).ColumnWidth = 29.144371396019366

# zh-cn: Status column widened, Error Detail column narrowed (now empty)
$ws_zhcn.Columns.Item(3).ColumnWidth = 29.144371396019366
$ws_zhcn.Columns.Item(16).ColumnWidth = 12.913719540550566

# de-de: Status column widened, Error Detail column narrowed (now empty)
$ws_dede.Columns.Item(3).ColumnWidth = 29.144371396019366
$ws_dede.Columns.Item(16).ColumnWidth = 12.913719540550566
